$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New content: "2 dices, at least one 4" block (rows 3-5, col I/J) ---
$ws.Range("I3").Value = "2 dices, at least one 4"

$ws.Range("I4").NumberFormat = "# ?/?"
$ws.Range("I4").Font.Bold = $true

$ws.Range("L9").NumberFormat = "d-mmm"
$ws.Range("L9").Font.Bold = $true
$ws.Range("L9").Value = "2-6,3-5,4-4,5-3,6-2"

$ws.Range("J5").Formula = "=11/36"
$ws.Range("J5").NumberFormat = "# ??/??"
$ws.Range("J5").Font.Bold = $true

$ws.Range("J10").Formula = "=1/9+1/8+1/17"

# --- New content: probability table (rows 17-20, cols J/K/L/M) ---
# Shared-string insertion order matches the original authoring order:
# J17, K17, K18, K20, then J18 last.
$ws.Range("J17").Value = "2,3,5,7"
$ws.Range("K17").Value = "p(a)"
$ws.Range("L17").Formula = '="4/10"'

$ws.Range("K18").Value = "p(b)"
$ws.Range("L18").NumberFormat = "d-mmm"
$ws.Range("L18").Font.Bold = $true
$ws.Range("L18").Formula = '="5/10"'
$ws.Range("M18").Formula = '="20/100 or 2/10"'

$ws.Range("K20").Value = "p(a and b)"
$ws.Range("L20").NumberFormat = "d-mmm"
$ws.Range("L20").Font.Bold = $true
$ws.Range("L20").Formula = '="3/10"'

$ws.Range("J18").Value = "1,2,3,4,5"

# --- View state ---
$ws.Range("M13").Select()
